# Bug fix astype or int values and load solved data
$wb = $excel.ActiveWorkbook

$wsScenarios = $wb.Worksheets.Item("scenarios")
$wsSub = $wb.Worksheets.Item("sub-scenarios")

# ---- "scenarios" sheet: rewrite B2:C2 and add rows 3-5 ----
$wsScenarios.Range("B2").Value = "Portfolio P0_2030"
$wsScenarios.Range("C2").Value = "input_P0_2030"

$wsScenarios.Range("A3").Value = 1
$wsScenarios.Range("B3").Value = "Portfolio P1_2030"
$wsScenarios.Range("C3").Value = "input_P1_2030"

$wsScenarios.Range("A4").Value = 2
$wsScenarios.Range("B4").Value = "Portfolio P2_2030"
$wsScenarios.Range("C4").Value = "input_P2_2030"

$wsScenarios.Range("A5").Value = 3
$wsScenarios.Range("B5").Value = "Portfolio P3_2030"
$wsScenarios.Range("C5").Value = "input_P3_2030"

# match the style used by row 2 for the newly added rows (A3:A5)
$wsScenarios.Range("A2").Copy()
$wsScenarios.Range("A3:A5").PasteSpecial(-4122)

# widen column C a bit (bestFit after the longer strings)
$wsScenarios.Columns.Item(3).ColumnWidth = 13

# selection / tab state: "scenarios" is no longer the active tab
$wsScenarios.Range("A9").Select()

# ---- "sub-scenarios" sheet: update C2:C5 ----
$wsSub.Range("C2").Value = 2017
$wsSub.Range("C3").Value = "Moderat"
$wsSub.Range("C4").Value = "Hoch"
$wsSub.Range("C5").Value = "Niedrig"

# "sub-scenarios" becomes the active/selected tab
$wsSub.Activate()
$wsSub.Range("C3").Select()
